# Generate Report for Handback
# Updates the handoff/handback timestamps for the d6c199f8-... file across
# the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the d6c199f8 row (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-13 08:57:36"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the d6c199f8 row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-13 08:57:29"
$wsZhCn.Range("K3").Value = "2016-08-13 08:57:56"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the d6c199f8 row (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-13 08:57:36"
$wsDeDe.Range("K3").Value = "2016-08-13 08:58:10"
